# Q3 Update - 2025
# Applies the quarterly refresh to the UN-COB worksheet:
#  1. Inserts a new data row for "Burkina Faso" (between Burundi and Central African Rep.)
#  2. Refreshes several statistic cells (N/O/P/T columns) for the affected countries
#  3. Updates the report/page hash stored in column B for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at 472 (shifts the existing 472-486 rows down to 473-487)
# ---------------------------------------------------------------------------
$ws.Rows("472").Insert()

# Populate the newly inserted row with the "Burkina Faso" record
$ws.Range("A472").Value = "1"
$ws.Range("B472").Value = "Ak5c2K"
$ws.Range("C472").Value = "1"
$ws.Range("D472").Value = "471"
$ws.Range("E472").Value = "2024"
$ws.Range("F472").Value = "23"
$ws.Range("G472").Value = "Burkina Faso"
$ws.Range("H472").Value = "BKF"
$ws.Range("I472").Value = "BFA"
$ws.Range("J472").Value = "40"
$ws.Range("K472").Value = "Congo"
$ws.Range("L472").Value = "COB"
$ws.Range("M472").Value = "COG"
$ws.Range("N472").Value = "0"
$ws.Range("O472").Value = "7"
$ws.Range("P472").Value = "0"
$ws.Range("Q472").Value = "0"
$ws.Range("R472").Value = "0"
$ws.Range("S472").Value = "0"
$ws.Range("T472").Value = "0"
$ws.Range("U472").Value = "-"
$ws.Range("V472").Value = "0"

# ---------------------------------------------------------------------------
# 2. Refresh statistic values for the rows whose figures changed
#    (row numbers below already reflect the post-insert layout)
# ---------------------------------------------------------------------------

# Angola
$ws.Range("N469").Value = "25"

# Burundi
$ws.Range("N471").Value = "93"
$ws.Range("T471").Value = "18"

# Central African Rep.
$ws.Range("N473").Value = "29975"
$ws.Range("O473").Value = "5246"
$ws.Range("P473").Value = "46"
$ws.Range("T473").Value = "2906"

# Chad
$ws.Range("O474").Value = "98"

# Cameroon
$ws.Range("N475").Value = "6"
$ws.Range("O475").Value = "73"

# Congo (self)
$ws.Range("T476").Value = "100"

# Dem. Rep. of the Congo
$ws.Range("N477").Value = "29005"
$ws.Range("O477").Value = "1441"
$ws.Range("T477").Value = "3019"

# Cote d'Ivoire
$ws.Range("O479").Value = "20"

# Liberia
$ws.Range("N480").Value = "8"

# Mauritania
$ws.Range("N481").Value = "11"
$ws.Range("O481").Value = "45"
$ws.Range("T481").Value = "5"

# Mali
$ws.Range("O482").Value = "6"

# Nigeria
$ws.Range("O483").Value = "16"
$ws.Range("T483").Value = "5"

# Rwanda
$ws.Range("N484").Value = "2812"
$ws.Range("O484").Value = "427"
$ws.Range("T484").Value = "7297"

# South Sudan
$ws.Range("N485").Value = "5"
$ws.Range("O485").Value = "11"

# Sudan
$ws.Range("N486").Value = "44"
$ws.Range("O486").Value = "249"

# Syrian Arab Rep.
$ws.Range("N487").Value = "19"

# ---------------------------------------------------------------------------
# 3. Refresh the report/page hash (column B) stored on every data row
# ---------------------------------------------------------------------------
$lastRow = $ws.Cells(1048576, 1).End(-4162).Row
$ws.Range("B2:B" + $lastRow).Value = "Ak5c2K"
